$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new price observation is inserted as the new row 28,
# pushing all the existing historical rows (old 28-43) down by one (new 29-44).
$ws.Rows("28:28").Insert()

$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(28, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(28, 4).Value = 45089
$ws.Cells.Item(28, 5).Value = 15
$ws.Cells.Item(28, 6).Value = 100112003
$ws.Cells.Item(28, 7).Value = "Ajo"
$ws.Cells.Item(28, 8).Value = "Chino"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 1350
$ws.Cells.Item(28, 11).Value = 17000
$ws.Cells.Item(28, 12).Value = 19000
$ws.Cells.Item(28, 13).Value = 17741
$ws.Cells.Item(28, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(28, 15).Value = "China"
$ws.Cells.Item(28, 16).Value = 1774
$ws.Cells.Item(28, 17).Value = 10
$ws.Cells.Item(28, 18).Value = "Hortaliza"
